$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data (columns B:AC) between each of these consecutive row pairs.
# Column A (the positional record counter) is left untouched on each row.
$swapPairs = @(
    @(112, 113),
    @(122, 123),
    @(250, 251),
    @(254, 255)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")

    $vals1 = $range1.Value()
    $vals2 = $range2.Value()

    $range1.Value = $vals2
    $range2.Value = $vals1
}

# Row 295 (match 7824497, Sportivo Ameliano vs Sportivo Luqueno) was removed
# from the source feed; every following row shifts up one position.
$ws.Rows(295).Delete()

# Column A is a plain positional counter (row number - 2); restore it for the
# rows that shifted up.
for ($r = 295; $r -le 298; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
